$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their literal text representation
# (values like "291.36" or "1.003" would otherwise be auto-parsed as numbers)
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "22.443.46"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "1.573.53"
$ws.Range("E3").Value = "  +0.01%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.03%  "

# Row 6
$ws.Range("D6").Value = "291.36"
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").Value = "0.3740"
$ws.Range("E7").Value = "  -0.74%  "

# Row 8
$ws.Range("D8").Value = "49.88"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  -0.76%  "

# Row 10
$ws.Range("D10").Value = "0.07557"
$ws.Range("E10").Value = "  -1.42%  "

# Row 11
$ws.Range("D11").Value = "1.137"
$ws.Range("E11").Value = "  -2.08%  "

# Row 12
$ws.Range("E12").Value = "  +0.00%  "

# Row 13
$ws.Range("E13").Value = "  +0.63%  "

# Row 14
$ws.Range("D14").Value = "5.996"
$ws.Range("E14").Value = "  -0.23%  "

# Row 15
$ws.Range("D15").Value = "6.942"
$ws.Range("E15").Value = "  -0.02%  "

# Row 16
$ws.Range("D16").Value = "1.568.37"
$ws.Range("E16").Value = "  -0.22%  "

# Row 17
$ws.Range("E17").Value = "  -1.03%  "

# Row 18
$ws.Range("D18").Value = "91.13"
$ws.Range("E18").Value = "  +0.89%  "

# Row 19
$ws.Range("D19").Value = "0.06736"
$ws.Range("E19").Value = "  -0.46%  "

# Row 20
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("D21").Value = "6.264"
$ws.Range("E21").Value = "  +0.58%  "

# Row 22
$ws.Range("E22").Value = "  -2.45%  "

# Row 23
$ws.Range("D23").Value = "12.15"
$ws.Range("E23").Value = "  +0.94%  "

# Row 24
$ws.Range("D24").Value = "22.448.88"
$ws.Range("E24").Value = "  +0.12%  "

# Row 25
$ws.Range("D25").Value = "2.331"
$ws.Range("E25").Value = "  -4.06%  "

# Row 26
$ws.Range("D26").Value = "2.605"
$ws.Range("E26").Value = "  -4.72%  "

# Row 27
$ws.Range("D27").Value = "20.15"
$ws.Range("E27").Value = "  -0.80%  "

# Row 28
$ws.Range("D28").Value = "148.60"
$ws.Range("E28").Value = "  +1.73%  "

# Row 29
$ws.Range("D29").Value = "5.003"
$ws.Range("E29").Value = "  -0.59%  "

# Row 30
$ws.Range("D30").Value = "125.73"
$ws.Range("E30").Value = "  -0.58%  "

# Row 31
$ws.Range("D31").Value = "1.743.77"
$ws.Range("E31").Value = "  -0.16%  "

# Row 32
$ws.Range("D32").Value = "1.053"
$ws.Range("E32").Value = "  +5.11%  "

# Row 33
$ws.Range("D33").Value = "6.117"
$ws.Range("E33").Value = "  -1.42%  "

# Row 34
$ws.Range("D34").Value = "1.981"
$ws.Range("E34").Value = "  -1.80%  "

# Row 35
$ws.Range("D35").Value = "9.819"
$ws.Range("E35").Value = "  -2.31%  "

# Row 36
$ws.Range("D36").Value = "0.08403"
$ws.Range("E36").Value = "  -1.75%  "

# Row 37
$ws.Range("E37").Value = "  +3.53%  "

# Row 38
$ws.Range("D38").Value = "0.02469"
$ws.Range("E38").Value = "  -3.17%  "

# Row 39
$ws.Range("D39").Value = "0.2289"
$ws.Range("E39").Value = "  -1.25%  "

# Row 40
$ws.Range("D40").Value = "0.06513"
$ws.Range("E40").Value = "  -1.27%  "

# Row 41
$ws.Range("D41").Value = "5.460"

# Row 42
$ws.Range("E42").Value = "  -2.57%  "

# Row 43
$ws.Range("D43").Value = "0.6254"
$ws.Range("E43").Value = "  -3.08%  "

# Row 44
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$ws.Range("D45").Value = "13.97"
$ws.Range("E45").Value = "  -1.08%  "

# Row 46
$ws.Range("D46").Value = "3.811"
$ws.Range("E46").Value = "  +0.31%  "

# Row 47
$ws.Range("D47").Value = "0.5818"
$ws.Range("E47").Value = "  -3.29%  "

# Row 50
$ws.Range("D50").Value = "1.222"
$ws.Range("E50").Value = "  -6.14%  "

# Row 51
$ws.Range("D51").Value = "0.07323"
$ws.Range("E51").Value = "  -0.07%  "

# Rows 48-49: Quant and NEARProtocol swap positions with updated values
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "129.48"
$ws.Range("E48").Value = "  +3.19%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.082"
$ws.Range("E49").Value = "  -0.28%  "

# Restore default (unstyled) cell style now that the text values are locked in
$dataRange.Style = "Normal"
